# Update Pais sheet: reorder Benin/Surinam/Sierra Leona and Guyana/Taiwan,
# refresh the "Datos actualizados" timestamp, and update the latest case
# counts for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados a ..." timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 00:23"

# --- Estados Unidos (row 4) ---
$ws.Cells.Item(4,2).Value = 4911245
$ws.Cells.Item(4,3).Value = 47329
$ws.Cells.Item(4,4).Value = 2475493
$ws.Cells.Item(4,5).Value = 2275700
$ws.Cells.Item(4,7).Value = 1124
$ws.Cells.Item(4,8).Value = 160052

# --- Brasil (row 5) ---
$ws.Cells.Item(5,4).Value = 1970767
$ws.Cells.Item(5,5).Value = 735335

# --- Colombia (row 13) ---
$ws.Cells.Item(13,2).Value = 334979
$ws.Cells.Item(13,3).Value = 7129
$ws.Cells.Item(13,4).Value = 180258
$ws.Cells.Item(13,5).Value = 143406
$ws.Cells.Item(13,7).Value = 298
$ws.Cells.Item(13,8).Value = 11315

# --- Egipto (row 29) ---
$ws.Cells.Item(29,2).Value = 94752
$ws.Cells.Item(29,3).Value = 112
$ws.Cells.Item(29,4).Value = 45569
$ws.Cells.Item(29,5).Value = 44271
$ws.Cells.Item(29,7).Value = 24
$ws.Cells.Item(29,8).Value = 4912

# --- Barein (row 52) ---
$ws.Cells.Item(52,2).Value = 42132
$ws.Cells.Item(52,3).Value = 297
$ws.Cells.Item(52,4).Value = 39335
$ws.Cells.Item(52,5).Value = 2646

# --- Suiza (row 58) ---
$ws.Cells.Item(58,4).Value = 31600
$ws.Cells.Item(58,5).Value = 2165

# --- Bulgaria (row 81) ---
$ws.Cells.Item(81,2).Value = 12414
$ws.Cells.Item(81,3).Value = 255
$ws.Cells.Item(81,4).Value = 6964
$ws.Cells.Item(81,5).Value = 5035
$ws.Cells.Item(81,7).Value = 11
$ws.Cells.Item(81,8).Value = 415

# --- Rows 133-135: Surinam / Sierra Leona / Benin reordered to
#     Benin / Surinam / Sierra Leona, with refreshed counts ---
$ws.Cells.Item(133,1).Value = "Benin"
$ws.Cells.Item(133,2).Value = 1914
$ws.Cells.Item(133,3).Value = 109
$ws.Cells.Item(133,4).Value = 1036
$ws.Cells.Item(133,5).Value = 840
$ws.Cells.Item(133,7).Value = 2
$ws.Cells.Item(133,8).Value = 38

$ws.Cells.Item(134,1).Value = "Surinam"
$ws.Cells.Item(134,2).Value = 1893
$ws.Cells.Item(134,3).Value = 0
$ws.Cells.Item(134,4).Value = 1227
$ws.Cells.Item(134,5).Value = 639
$ws.Cells.Item(134,7).Value = 0
$ws.Cells.Item(134,8).Value = 27

$ws.Cells.Item(135,1).Value = "Sierra Leona"
$ws.Cells.Item(135,2).Value = 1855
$ws.Cells.Item(135,3).Value = 7
$ws.Cells.Item(135,4).Value = 1397
$ws.Cells.Item(135,5).Value = 391
$ws.Cells.Item(135,7).Value = 0
$ws.Cells.Item(135,8).Value = 67

# --- Santo Tome y Principe (row 155) ---
$ws.Cells.Item(155,2).Value = 875
$ws.Cells.Item(155,3).Value = 1
$ws.Cells.Item(155,4).Value = 794
$ws.Cells.Item(155,5).Value = 66

# --- Rows 165-166: Taiwan / Guyana reordered to Guyana / Taiwan,
#     with refreshed counts ---
$ws.Cells.Item(165,1).Value = "Guyana"
$ws.Cells.Item(165,2).Value = 497
$ws.Cells.Item(165,3).Value = 23
$ws.Cells.Item(165,4).Value = 186
$ws.Cells.Item(165,5).Value = 289
$ws.Cells.Item(165,7).Value = 1
$ws.Cells.Item(165,8).Value = 22

$ws.Cells.Item(166,1).Value = "Taiwan"
$ws.Cells.Item(166,2).Value = 476
$ws.Cells.Item(166,3).Value = 1
$ws.Cells.Item(166,4).Value = 441
$ws.Cells.Item(166,5).Value = 28
$ws.Cells.Item(166,7).Value = 0
$ws.Cells.Item(166,8).Value = 7
